$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for the season record columns, copying the
# formatting (bold font, border, centered alignment) from the existing
# header row so the new cells share the same style as AC1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Populate the season record (Wins/Losses/Ties) for every player row.
# Every row shares the same team record.
$ws.Range("AD2:AD50").Value = 69
$ws.Range("AE2:AE50").Value = 93
$ws.Range("AF2:AF50").Value = 0

Write-Host "done"
